$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2..92) holds a date serial that needs to move from
# 2023-09-02 (45171) to 2023-09-03 (45172).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value = 45172
    }
}
